$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.175199747085571
$ws.Range("B1").Value = 2.169065952301025
$ws.Range("C1").Value = 3.201991081237793
$ws.Range("D1").Value = 3.738319396972656
$ws.Range("E1").Value = 1.235240578651428
